$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its text representation instead of
# being auto-coerced to a number by COM type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.197.52"
$ws.Range("E2").Value = "  +3.90%  "
$ws.Range("D3").Value = "1.786.49"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "339.29"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "0.3829"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("D8").Value = "0.3441"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "47.08"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "1.152"
$ws.Range("E10").Value = "  -2.52%  "
$ws.Range("D11").Value = "0.07404"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "23.20"
$ws.Range("E12").Value = "  +8.08%  "
$ws.Range("D13").Value = "1.002"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "6.464"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "7.413"
$ws.Range("E15").Value = "  +4.86%  "
$ws.Range("D16").Value = "1.780.14"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "0.00001078"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "0.06692"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "82.36"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "17.52"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "6.477"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "28.234.86"
$ws.Range("E23").Value = "  +4.04%  "
$ws.Range("D24").Value = "12.09"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").Value = "2.368"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "1.451"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "20.74"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "2.424"
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("D29").Value = "154.36"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "136.38"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").Value = "1.985.71"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "6.145"
$ws.Range("E32").Value = "  +3.31%  "
$ws.Range("D33").Value = "3.959"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "0.08930"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("D35").Value = "12.80"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").Value = "0.02424"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").Value = "0.6866"
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").Value = "5.342"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").Value = "0.06382"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").Value = "0.2176"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").Value = "1.248"
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("D42").Value = "1.499"
$ws.Range("E42").Value = "  -7.12%  "
$ws.Range("D43").Value = "8.291"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "14.26"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "0.6313"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").Value = "3.876"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "133.48"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").Value = "0.07527"
$ws.Range("E50").Value = "  +6.02%  "
$ws.Range("D51").Value = "1.212"
$ws.Range("E51").Value = "  +6.75%  "

# Restore the default (un-styled) cell style now that the values are
# committed as text, matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
